$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected NN values for rows 2-9 (columns B, C, D, E, F, H)
# Column G (Площадь / area) is unchanged.

$data = @{
    2 = @{ B = 1245.60595703125;   C = 0.9344;                 D = 0.9050999879837036; E = 1.166499972343445;  F = 0.8461999893188477; H = 0.8451 }
    3 = @{ B = 1133.99609375;      C = 0.9072;                 D = 0.8989;              E = 1.233299970626831;  F = 0.8403000235557556; H = 0.7897999999999999 }
    4 = @{ B = 762.0167846679688;  C = 0.8902;                 D = 0.8903;              E = 0.9365000128746033; F = 0.8271999955177307; H = 0.7144 }
    5 = @{ B = 848.780029296875;   C = 0.8953;                 D = 0.8949;              E = 1.111899971961975;  F = 0.8492000102996826; H = 0.7547 }
    6 = @{ B = 1184.694091796875;  C = 0.9373;                 D = 0.9281;              E = 1.168300032615662;  F = 0.8271999955177307; H = 1.0485 }
    7 = @{ B = 927.5018920898438;  C = 0.9340000000000001;     D = 0.9291999936103821;  E = 1.074300050735474;  F = 0.885200023651123;  H = 1.0589 }
    8 = @{ B = 1036.651733398438;  C = 0.9288999999999999;     D = 0.9249000000000001;  E = 1.100200057029724;  F = 0.8867999911308289; H = 1.021 }
    9 = @{ B = 7139.24658203125;   C = 0.92;                   D = 0.9121;              E = 1.233299970626831;  F = 0.8271999955177307; H = 6.232399999999999 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("H$row").Value = $vals.H
}

$wb.Save()
